$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 195, shifting existing rows 195:200 down to 196:201
$ws.Rows.Item(195).Insert()

# Populate the new row 195 with the new weekly price record
$ws.Range("A195").Value2 = 8
$ws.Range("B195").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C195").Value2 = "Coquimbo"
$ws.Range("D195").Value2 = 44568
$ws.Range("D195").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E195").Value2 = 4
$ws.Range("F195").Value2 = 100112003
$ws.Range("G195").Value2 = "Ajo"
$ws.Range("H195").Value2 = "Chino"
$ws.Range("I195").Value2 = "Primera"
$ws.Range("J195").Value2 = 600
$ws.Range("K195").Value2 = 18000
$ws.Range("L195").Value2 = 19000
$ws.Range("M195").Value2 = 18500
$ws.Range("N195").Value2 = "$/caja 10 kilos"
$ws.Range("O195").Value2 = "China"
$ws.Range("P195").Value2 = 1850
$ws.Range("Q195").Value2 = 10
$ws.Range("R195").Value2 = "Hortaliza"
